$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style of the existing H1 header cell
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J21
$iValues = @(1,1,1,1,1,1,1,7,5,7,8,1,3,6,7,8,1,1,1,1)
$jValues = @(4,4,4,5,7,5,6,8,6,8,9,4,7,7,7,9,4,4,4,2)

for ($r = 0; $r -lt 20; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
